# Update Betfair Back/Lay odds and liquidity figures on the active sheet
# (rows 2-12 of the "Jogos do Dia" odds table) to reflect the latest snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2
$ws.Range("H2").Value = 4.3
$ws.Range("I2").Value = 5.2
$ws.Range("J2").Value = 2.96
$ws.Range("K2").Value = 3.7
$ws.Range("L2").Value = 1.5
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 2.42
$ws.Range("O2").Value = 1.57
$ws.Range("P2").Value = 1.47
$ws.Range("Q2").Value = 2.68
$ws.Range("R2").Value = 1.16
$ws.Range("T2").Value = 2.24
$ws.Range("U2").Value = 1.65
$ws.Range("Z2").Value = 1000
$ws.Range("AB2").Value = 14
$ws.Range("AK2").Value = 1000
$ws.Range("F3").Value = 3.05
$ws.Range("H3").Value = 2.72
$ws.Range("J3").Value = 2.74
$ws.Range("N3").Value = 2.22
$ws.Range("O3").Value = 1.66
$ws.Range("P3").Value = 1.4
$ws.Range("U3").Value = 1.64
$ws.Range("V3").Value = 1.48
$ws.Range("AB3").Value = 19.5
$ws.Range("AC3").Value = 25
$ws.Range("G4").Value = 1.16
$ws.Range("H4").Value = 21
$ws.Range("J4").Value = 10
$ws.Range("K4").Value = 13
$ws.Range("P4").Value = 3.85
$ws.Range("S4").Value = 1.73
$ws.Range("T4").Value = 1.95
$ws.Range("W4").Value = 7.2
$ws.Range("AD4").Value = 85
$ws.Range("AK4").Value = 16
$ws.Range("AN4").Value = 2.66
$ws.Range("F5").Value = 1.85
$ws.Range("G5").Value = 2
$ws.Range("L5").Value = 1.56
$ws.Range("N5").Value = 2.52
$ws.Range("O5").Value = 1.54
$ws.Range("P5").Value = 1.51
$ws.Range("Q5").Value = 2.58
$ws.Range("R5").Value = 1.18
$ws.Range("T5").Value = 2.24
$ws.Range("U5").Value = 1.66
$ws.Range("V5").Value = 1.2
$ws.Range("W5").Value = 2
$ws.Range("Y5").Value = 14
$ws.Range("AB5").Value = 6.4
$ws.Range("G6").Value = 2.64
$ws.Range("Q6").Value = 2.44
$ws.Range("AF6").Value = 15
$ws.Range("AH6").Value = 20
$ws.Range("AO6").Value = 60
$ws.Range("H7").Value = 5.2
$ws.Range("I7").Value = 5.8
$ws.Range("K7").Value = 5.3
$ws.Range("P7").Value = 3.3
$ws.Range("Q7").Value = 1.35
$ws.Range("R7").Value = 1.95
$ws.Range("T7").Value = 1.5
$ws.Range("U7").Value = 2.72
$ws.Range("X7").Value = 42
$ws.Range("AO7").Value = 34
$ws.Range("L8").Value = 1.33
$ws.Range("N8").Value = 3.7
$ws.Range("Q8").Value = 1.92
$ws.Range("S8").Value = 3.35
$ws.Range("V8").Value = 2.16
$ws.Range("W8").Value = 1.21
$ws.Range("X8").Value = 27
$ws.Range("Y8").Value = 10
$ws.Range("AA8").Value = 900
$ws.Range("AB8").Value = 980
$ws.Range("AF8").Value = 95
$ws.Range("AG8").Value = 980
$ws.Range("G9").Value = 1.72
$ws.Range("H9").Value = 7
$ws.Range("I9").Value = 8.4
$ws.Range("J9").Value = 3.5
$ws.Range("K9").Value = 3.9
$ws.Range("L9").Value = 1.53
$ws.Range("M9").Value = 1.11
$ws.Range("N9").Value = 2.62
$ws.Range("O9").Value = 1.51
$ws.Range("P9").Value = 1.55
$ws.Range("Q9").Value = 2.48
$ws.Range("R9").Value = 1.19
$ws.Range("S9").Value = 5
$ws.Range("T9").Value = 2.34
$ws.Range("U9").Value = 1.6
$ws.Range("V9").Value = 1.13
$ws.Range("W9").Value = 2.38
$ws.Range("AB9").Value = 6
$ws.Range("AF9").Value = 8.6
$ws.Range("AG9").Value = 11
$ws.Range("AK9").Value = 25
$ws.Range("AL9").Value = 160
$ws.Range("AN9").Value = 17.5
$ws.Range("I10").Value = 10.5
$ws.Range("N10").Value = 3.5
$ws.Range("O10").Value = 1.38
$ws.Range("T10").Value = 2.4
$ws.Range("AA10").Value = 480
$ws.Range("AB10").Value = 6.4
$ws.Range("AC10").Value = 11
$ws.Range("AF10").Value = 7.2
$ws.Range("AK10").Value = 17
$ws.Range("AL10").Value = 50
$ws.Range("AN10").Value = 8.6
$ws.Range("AO10").Value = 380
$ws.Range("H11").Value = 6.6
$ws.Range("I11").Value = 6.8
$ws.Range("N11").Value = 4.5
$ws.Range("O11").Value = 1.26
$ws.Range("P11").Value = 2.2
$ws.Range("Q11").Value = 1.8
$ws.Range("R11").Value = 1.47
$ws.Range("S11").Value = 3.05
$ws.Range("T11").Value = 1.89
$ws.Range("U11").Value = 2.08
$ws.Range("W11").Value = 2.68
$ws.Range("Y11").Value = 24
$ws.Range("Z11").Value = 60
$ws.Range("AB11").Value = 9
$ws.Range("AC11").Value = 9.6
$ws.Range("AF11").Value = 9.199999999999999
$ws.Range("AG11").Value = 9.4
$ws.Range("AH11").Value = 21
$ws.Range("AL11").Value = 30
$ws.Range("AM11").Value = 95
$ws.Range("F12").Value = 2.18
$ws.Range("J12").Value = 3.15
$ws.Range("N12").Value = 2.86
$ws.Range("S12").Value = 4.5
$ws.Range("U12").Value = 1.89
$ws.Range("V12").Value = 1.31
$ws.Range("AH12").Value = 21
